# Apply the "Podredba v osnovnata rabota" edit described by the commit:
#  - mark the title paragraph run properties with en-US language
#  - split the title run (fixing "приложимостн" -> "приложимост") into three runs
#  - tag the picture run with noProof/eastAsia language info
#  - append the new body paragraphs (results / effectiveness / applicability / trailing blank)
$d = $word.ActiveDocument
$target = $d.Content

$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships">
        <w:body>
<w:p w:rsidR="003A1C3C" w:rsidRPr="003A1C3C" w:rsidRDefault="003A1C3C" w:rsidP="003A1C3C"><w:pPr><w:jc w:val="right"/><w:rPr><w:b/><w:sz w:val="36"/><w:szCs w:val="36"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r w:rsidRPr="003A1C3C"><w:rPr><w:noProof/><w:lang w:eastAsia="bg-BG"/></w:rPr><w:drawing><wp:anchor distT="0" distB="0" distL="114300" distR="114300" simplePos="0" relativeHeight="251658240" behindDoc="0" locked="0" layoutInCell="1" allowOverlap="1"><wp:simplePos x="0" y="0"/><wp:positionH relativeFrom="column"><wp:posOffset>22860</wp:posOffset></wp:positionH><wp:positionV relativeFrom="paragraph"><wp:posOffset>-368300</wp:posOffset></wp:positionV><wp:extent cx="5775325" cy="1860550"/><wp:effectExtent l="19050" t="0" r="0" b="0"/><wp:wrapTopAndBottom/><wp:docPr id="69" name="Picture 4" descr="glava5.JPG"/><wp:cNvGraphicFramePr><a:graphicFrameLocks xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" noChangeAspect="1"/></wp:cNvGraphicFramePr><a:graphic xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main"><a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:pic xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:nvPicPr><pic:cNvPr id="0" name="glava5.JPG"/><pic:cNvPicPr/></pic:nvPicPr><pic:blipFill><a:blip r:embed="rId4" cstate="print"/><a:srcRect l="3871" r="2841"/><a:stretch><a:fillRect/></a:stretch></pic:blipFill><pic:spPr><a:xfrm><a:off x="0" y="0"/><a:ext cx="5775325" cy="1860550"/></a:xfrm><a:prstGeom prst="rect"><a:avLst/></a:prstGeom></pic:spPr></pic:pic></a:graphicData></a:graphic></wp:anchor></w:drawing></w:r><w:r w:rsidRPr="003A1C3C"><w:rPr><w:b/><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr><w:t>Оценка на получениете резултати, техни</w:t></w:r><w:r><w:rPr><w:b/><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr><w:t>ческа ефективност и приложимост</w:t></w:r><w:r><w:rPr><w:b/><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr><w:t xml:space="preserve"> на дипломната работа</w:t></w:r></w:p><w:p><w:pPr><w:jc w:val="both"/><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:jc w:val="both"/><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>Резулататите от експериментите са очаквани и задоволителни. Те показват зависимостта между отделните параметри и влиянието им върху полета и стабилизацията на трикоптера.</w:t></w:r></w:p><w:p><w:pPr><w:jc w:val="both"/><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>Ефективността на трикоптера с каскадната управляваща система е достатъчно добра, за да позволява изпълнението на сложни мисии с трикоптера.</w:t></w:r></w:p><w:p><w:pPr><w:jc w:val="both"/><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve">Приложимост дипломната работа може да намери като въвеждане и подробно описание на летателен апарат и начините му за управление чрез Ардуино. Дипломната работа може да служи и като въвведение в света на ардуино.   </w:t></w:r></w:p><w:p><w:pPr><w:jc w:val="both"/><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve">  </w:t></w:r></w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

[void]$target.InsertXML($xml)
